$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('070028236358-54', 'ASML Netherlands B.V. (ASML)'),
    @('583393138546-76', 'Applied Materials (AMAT)'),
    @('97921387213-35', 'Interuniversitair Micro-Electronica Centrum (IMEC)'),
    @('866541247974-60', 'NVIDIA Corporation'),
    @('00358442856-45', 'QUALCOMM Inc. (QCOM (NASDAQ))'),
    @('389630949944-21', 'ARDIAN France'),
    @('236569017571-35', 'Association Nationale des Sociétés par Actions (ANSA)'),
    @('7459401905-60', 'Intel Corporation'),
    @('809121549545-61', 'Taiwan Semiconductor Manufacturing Company Ltd'),
    @('236869539281-13', 'The LEGO Group'),
    @('70843878143-13', 'STMicroelectronics (ST)'),
    @('012717123539-58', 'European Chamber of Commerce Taiwan (ECCT)'),
    @('494613715191-85', 'Cisco Systems Inc. (Cisco)'),
    @('524976815651-12', 'Centre for European Reform (CER)'),
    @('194094423131-36', 'ZF Friedrichshafen AG'),
    @('547225616323-13', 'Tesla, Inc. (Tesla)'),
    @('89093924456-06', 'Vereinigung der österreichischen Industrie - Industriellenvereinigung (IV)'),
    @('27762251795-15', 'Confederazione Generale dell''Industria Italiana (Confindustria)'),
    @('142903829084-58', 'Stiftung Neue Verantwortung e.V. (SNV)'),
    @('001968511387-52', 'European Internet Forum (EIF)'),
    @('1771817758-48', 'Bundesverband der Deutschen Industrie e.V. (BDI)'),
    @('475891045627-85', 'AT & S AUSTRIA TECHNOLOGIE & SYSTEMTECHNIK AKTIENGESELLSCHAFT'),
    @('336735239719-29', 'Kalray SA'),
    @('986044541551-20', 'Stellantis'),
    @('8999533555-91', 'Robert Bosch GmbH'),
    @('04201463642-88', 'APPLiA (Home Appliance Europe) (APPLiA)'),
    @('21856815315-64', 'EUROSMART'),
    @('35124436517-75', 'Fraunhofer-Gesellschaft (FhG)'),
    @('898711618618-45', 'NXP Semiconductors Netherlands B.V. (NXP)'),
    @('588327811384-96', 'Apple Inc.'),
    @('721066539792-75', 'GLOBALFOUNDRIES Management Services Limited Liability Company & Co. KG (GF Dresden)'),
    @('7721359944-96', 'IBM Corporation (IBM)'),
    @('40471017282-57', 'Samsung Electronics Europe'),
    @('266271142555-68', 'Taiwan Semiconductor Manufacturing Company Ltd (TSMC)'),
    @('3305029916-47', 'Scania AB (publ) (Scania)'),
    @('10751968675-85', 'Infineon Technologies AG'),
    @('532086816465-17', 'Arm Holdings'),
    @('789826549861-77', 'Synopsys Inc.'),
    @('22400601191-42', 'Deutscher Industrie- und Handelskammertag e.V. (DIHK)'),
    @('764184537594-67', 'SAFRAN'),
    @('174684235534-60', 'Qualcomm Communications SARL'),
    @('64270747023-20', 'DIGITALEUROPE (DE)'),
    @('20210641335-88', 'ORGALIME – The European Technology Industries (ORGALIME)'),
    @('5749958415-41', 'Confederation of Danish Industry (DI)'),
    @('865956522785-09', 'Fraunhofer Austria Reserach GmbH (FhA)'),
    @('248908039789-74', 'Melexis NV'),
    @('367850645092-18', 'Nordic Semiconductor ASA'),
    @('030658245107-49', 'Alcyon Photonics SL'),
    @('041592911733-05', 'Centre on Regulation in Europe (CERRE)'),
    @('71635694112-37', 'FEANI (FEANI)'),
    @('147656111703-83', 'European Forum of technical and Vocational Education and Training (EfVET)'),
    @('358284014848-82', 'industriAll European Trade Union (industriAll)'),
    @('52774696782-43', 'Commissariat à l''Energie atomique et aux énergies alternatives (CEA)'),
    @('817789336491-32', 'SOITEC (SOITEC)'),
    @('402302029423-14', 'SEMI Europe (SEMI)'),
    @('243961945295-86', 'SOITEC'),
)

$numRows = $data.Count
$startRow = 2
$endRow = $startRow + $numRows - 1

for ($i = 0; $i -lt $numRows; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
